# Auto-generated edit script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.579.07"
$ws.Range("D3").Value = "1.665.59"
$ws.Range("E3").Value = "  -3.52%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.51"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.514"
$ws.Range("E6").Value = "  -1.83%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.61"
$ws.Range("E8").Value = "  -2.06%  "
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("E10").Value = "  -1.61%  "
$ws.Range("E11").Value = "  -2.08%  "
$ws.Range("D12").Value = "1.900.68"
$ws.Range("E12").Value = "  -3.55%  "
$ws.Range("D13").Value = "1.656.03"
$ws.Range("E13").Value = "  -4.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.15"
$ws.Range("E14").Value = "  -2.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.560"
$ws.Range("E15").Value = "  -0.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.16"
$ws.Range("E16").Value = "  -2.07%  "
$ws.Range("D17").Value = "27.595.22"
$ws.Range("E17").Value = "  -1.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "242.25"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").Value = "0.0₃0730"
$ws.Range("E19").Value = "  -3.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.57"
$ws.Range("E20").Value = "  -3.82%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.49"
$ws.Range("E22").Value = "  -2.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.29"
$ws.Range("E23").Value = "  -4.52%  "
$ws.Range("E24").Value = "  -4.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.47"
$ws.Range("E25").Value = "  -1.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.19"
$ws.Range("E26").Value = "  -4.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.38"
$ws.Range("E27").Value = "  -2.17%  "
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.112"
$ws.Range("E28").Value = "  -2.16%  "
$ws.Range("B29").Value = "BinanceUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  +3.53%  "
$ws.Range("E31").Value = "  -1.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.35"
$ws.Range("E32").Value = "  -2.61%  "
$ws.Range("D33").Value = "1.476.87"
$ws.Range("E33").Value = "  -1.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.11"
$ws.Range("E34").Value = "  -4.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.57"
$ws.Range("E35").Value = "  -5.24%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.933"
$ws.Range("E36").Value = "  -1.94%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.38"
$ws.Range("E37").Value = "  -1.13%  "
$ws.Range("E38").Value = "  -5.33%  "
$ws.Range("E39").Value = "  -1.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "69.42"
$ws.Range("E40").Value = "  -2.00%  "
$ws.Range("E41").Value = "  -5.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("E43").Value = "  -7.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.21"
$ws.Range("E44").Value = "  -3.99%  "
$ws.Range("D45").Value = "1.808.87"
$ws.Range("E45").Value = "  -3.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.789"
$ws.Range("E46").Value = "  -0.88%  "
$ws.Range("E47").Value = "  -2.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "89.41"
$ws.Range("E48").Value = "  -1.87%  "
$ws.Range("E49").Value = "  -4.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.103"
$ws.Range("E50").Value = "  -2.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.90"
$ws.Range("E51").Value = "  -3.64%  "
